# fix: poster headline centering
#
# The "Master-Thesis" sub-headline text box on the poster slide was off
# center. Shift it left by 89284 EMU (1228725 -> 1139441) while keeping
# its vertical position and size unchanged.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$headline = $s.Shapes.Item("Text Box 126")

# EMU -> point conversion (1 pt = 12700 EMU) so the saved OOXML round-trips
# back to the exact target offset of 1139441 EMU.
$headline.Left = 1139441 / 12700
